# Applies the "Marksheet" recomputation edit to the quiz result workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary rows (10-12): recompute Right/Wrong/NotAttempt/Max, Marking and Total ---

# Row 10 "No." - counts
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 19
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 28

# Row 11 "Marking" - marks per right/wrong answer
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0

# Row 12 "Total" - totals and final score
$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 76
$ws.Range("C12").Value = -3
$ws.Range("E12").Value = "73/112"

# --- Answer key table (rows 15-40) ---

# Drop the third Student Ans / Correct Ans column group entirely (columns G:H)
$ws.Range("G15:H40").Clear()

# Fill in the "Student Ans" column (A) for the first question group, color-coded
# by correctness (green = correctStyle, red = incorrectStyle). Rows left out below
# (17, 24, 26, 29, 31, 35) were not attempted and stay blank.
$ws.Range("A16").Value = "Option A"
$ws.Range("A16").Style = "correctStyle"
$ws.Range("A18").Value = "Option B"
$ws.Range("A18").Style = "correctStyle"
$ws.Range("A19").Value = "Option C"
$ws.Range("A19").Style = "correctStyle"
$ws.Range("A20").Value = "Option B"
$ws.Range("A20").Style = "correctStyle"
$ws.Range("A21").Value = "Option C"
$ws.Range("A21").Style = "correctStyle"
$ws.Range("A22").Value = "Option D"
$ws.Range("A22").Style = "correctStyle"
$ws.Range("A23").Value = "Option D"
$ws.Range("A23").Style = "correctStyle"
$ws.Range("A25").Value = "Option A"
$ws.Range("A25").Style = "correctStyle"
$ws.Range("A27").Value = "Option A"
$ws.Range("A27").Style = "correctStyle"
$ws.Range("A28").Value = "Option D"
$ws.Range("A28").Style = "correctStyle"
$ws.Range("A30").Value = "Option B"
$ws.Range("A30").Style = "correctStyle"
$ws.Range("A32").Value = "Option C"
$ws.Range("A32").Style = "correctStyle"
$ws.Range("A33").Value = "Option D"
$ws.Range("A33").Style = "correctStyle"
$ws.Range("A34").Value = "Option B"
$ws.Range("A34").Style = "correctStyle"
$ws.Range("A36").Value = "Option C"
$ws.Range("A36").Style = "incorrectStyle"
$ws.Range("A37").Value = "Option C"
$ws.Range("A37").Style = "incorrectStyle"
$ws.Range("A38").Value = "Option A"
$ws.Range("A38").Style = "correctStyle"
$ws.Range("A39").Value = "Option D"
$ws.Range("A39").Style = "correctStyle"
$ws.Range("A40").Value = "Option D"
$ws.Range("A40").Style = "correctStyle"

# The second question group (columns D/E) only had 3 questions; fill in the
# "Student Ans" column (D) for those and drop the rest of the group (rows 19-40).
$ws.Range("D16").Value = "Option A"
$ws.Range("D16").Style = "correctStyle"
$ws.Range("D17").Value = "Option C"
$ws.Range("D17").Style = "correctStyle"
$ws.Range("D18").Value = "Option B"
$ws.Range("D18").Style = "incorrectStyle"
$ws.Range("D19:E40").Clear()
